# Add new test case for Profile settings.
# Current qa.test1@cvhcare.com password is Password2!
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FirstSet")
$ws.Range("B2").Value = "Password2!"
